$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update patient row 2 with new random data
$ws.Range("A2").Value = "test25"
$ws.Range("C2").Value = 1111111125

# Update patient row 3 with new random data
$ws.Range("A3").Value = "test26"
$ws.Range("C3").Value = 1111111126

# Update the selection to reflect the edited cell
$ws.Range("D2").Select()
